$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient log update (mohamad): vivalnk patch/charger #C600022 and
#     #C700136 are back at Biobank, not with Anish anymore.
$ws.Range("D4").Value = "Biobank"
$ws.Range("D13").Value = "Biobank"
$ws.Range("D14").Value = "Biobank"
$ws.Range("D18").Value = "Biobank"
$ws.Range("D19").Value = "Biobank"

# --- Updated patch logistics: newly logged vivalnk patches + chargers ---
$serials = @("C700138", "C700149", "C700157", "C700146", "C700148", "C700204", "C700205", "C700126", "C700127")

$row = 20
foreach ($serial in $serials) {
    $ws.Cells.Item($row, 1).Value = "patch_vivalnk"
    $ws.Cells.Item($row, 2).Value = $serial
    $ws.Cells.Item($row, 3).Value = "-"
    $ws.Cells.Item($row, 4).Value = "Biobank"
    $row++
}

foreach ($serial in $serials) {
    $ws.Cells.Item($row, 1).Value = "charger_vivalnk"
    $ws.Cells.Item($row, 2).Value = $serial
    $ws.Cells.Item($row, 3).Value = "-"
    $ws.Cells.Item($row, 4).Value = "Biobank"
    $row++
}

# Match the formatting used by the rest of the table (16pt font, 21pt rows).
$ws.Range("A20:D37").Font.Size = 16
$ws.Range("A20:D37").RowHeight = 21

# --- Restore selection to match the post-edit workbook state ---
$ws.Range("F22").Select()
